$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 is the J7 programming-header connector. The BOM is being revised
# to use a larger (.100" pitch) 6-position programming header instead of
# the old .050" pitch header.
$ws.Range("E21").Value = "609-3272-ND"
$ws.Range("F21").Value = "68001-106HLF"
$ws.Range("D21").Value = "CONN HEADER 6POS .100 STR 30AU"
$ws.Range("C21").Value = "MA06-1"
$ws.Range("H21").Value = 0.46

# Update the view state to reflect the scrolled/selected position after
# the edit.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C22").Select()
